# Finalizando script de cálculo de métricas
# Updates the "Changed files" list strings (column F) to drop the stray
# leading space after commas (and reorder one list), and recomputes the
# Precision/Recall/F2 (+Deps variants) metrics for rows 7, 12, 13 and 14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F: "TestIWithDeps" changed-files lists, reformatted ---
$ws.Range("F2").Value = "['app/controllers/home_controller.rb', 'app/models/question.rb', 'app/views/home/admin.html.haml', 'app/views/shared/_highcharts_header.html.haml']"
$ws.Range("F3").Value = "['app/models/choice.rb', 'app/models/earl.rb', 'app/models/item.rb', 'app/models/question.rb', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/views/abingo_dashboard/index.html.haml']"
$ws.Range("F4").Value = "['app/models/choice.rb', 'app/models/density.rb', 'app/models/earl.rb', 'app/models/prompt.rb', 'app/models/question.rb', 'app/models/session.rb', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/views/abingo_dashboard/index.html.haml']"
$ws.Range("F5").Value = "['app/models/choice.rb', 'app/models/density.rb', 'app/models/earl.rb', 'app/models/prompt.rb', 'app/models/question.rb', 'app/models/session.rb', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/views/abingo_dashboard/index.html.haml']"
$ws.Range("F6").Value = "['app/views/abingo_dashboard/index.html.haml', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/models/question.rb']"
$ws.Range("F7").Value = "['app/controllers/choices_controller.rb', 'app/models/choice.rb', 'app/models/earl.rb', 'app/models/prompt.rb', 'app/models/question.rb', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/views/abingo_dashboard/index.html.haml', 'app/views/choices/show.html.haml', 'app/views/shared/_header_vote.html.haml']"
$ws.Range("F8").Value = "['app/controllers/choices_controller.rb', 'app/controllers/home_controller.rb', 'app/controllers/questions_controller.rb', 'app/models/choice.rb', 'app/models/earl.rb', 'app/models/item.rb', 'app/models/prompt.rb', 'app/models/question.rb', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/views/abingo_dashboard/index.html.haml', 'app/views/choices/show.html.haml', 'app/views/home/about.html.haml', 'app/views/home/admin.html.haml', 'app/views/home/index.html.haml', 'app/views/home/privacy.html.haml', 'app/views/questions/_idea.html.haml', 'app/views/questions/about.html.haml', 'app/views/questions/admin.html.haml', 'app/views/questions/new.html.haml', 'app/views/questions/results.html.haml', 'app/views/questions/voter_map.html.erb', 'app/views/questions/word_cloud.html.erb', 'app/views/shared/_google_jsapi.html.haml', 'app/views/shared/_header_vote.html.haml', 'app/views/shared/_highcharts_header.html.haml']"
$ws.Range("F9").Value = "['app/models/question.rb', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/views/abingo_dashboard/index.html.haml']"
$ws.Range("F10").Value = "['app/models/choice.rb', 'app/models/prompt.rb', 'app/models/question.rb', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/views/abingo_dashboard/index.html.haml']"
$ws.Range("F11").Value = "['app/models/question.rb', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/views/abingo_dashboard/index.html.haml']"
$ws.Range("F12").Value = "['app/controllers/choices_controller.rb', 'app/controllers/home_controller.rb', 'app/controllers/questions_controller.rb', 'app/models/choice.rb', 'app/models/earl.rb', 'app/models/item.rb', 'app/models/prompt.rb', 'app/models/question.rb', 'app/models/user.rb', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/views/abingo_dashboard/index.html.haml', 'app/views/choices/show.html.haml', 'app/views/home/about.html.haml', 'app/views/home/admin.html.haml', 'app/views/home/index.html.haml', 'app/views/home/privacy.html.haml', 'app/views/questions/_idea.html.haml', 'app/views/questions/about.html.haml', 'app/views/questions/admin.html.haml', 'app/views/questions/new.html.haml', 'app/views/questions/results.html.haml', 'app/views/questions/voter_map.html.erb', 'app/views/questions/word_cloud.html.erb', 'app/views/shared/_google_jsapi.html.haml', 'app/views/shared/_header_vote.html.haml', 'app/views/shared/_highcharts_header.html.haml']"
$ws.Range("F13").Value = "['app/controllers/choices_controller.rb', 'app/models/choice.rb', 'app/models/earl.rb', 'app/models/prompt.rb', 'app/models/question.rb', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/views/abingo_dashboard/index.html.haml', 'app/views/choices/show.html.haml', 'app/views/shared/_header_vote.html.haml']"
$ws.Range("F14").Value = "['app/models/choice.rb', 'app/models/earl.rb', 'app/models/question.rb', 'app/views/abingo_dashboard/_experiment_row.html.haml', 'app/views/abingo_dashboard/index.html.haml']"

# --- Recomputed metrics for row 7 ---
$ws.Range("G7").Value = 0.1111111111111111
$ws.Range("H7").Value = 0.0625
$ws.Range("I7").Value = 0.06849315068493152
$ws.Range("J7").Value = 0.1111111111111111
$ws.Range("K7").Value = 0.0625
$ws.Range("L7").Value = 0.06849315068493152

# --- Recomputed metrics for row 12 ---
$ws.Range("G12").Value = 0.2692307692307692
$ws.Range("H12").Value = 0.1707317073170732
$ws.Range("I12").Value = 0.1842105263157895
$ws.Range("J12").Value = 0.2692307692307692
$ws.Range("K12").Value = 0.1707317073170732
$ws.Range("L12").Value = 0.1842105263157895

# --- Recomputed metrics for row 13 ---
$ws.Range("G13").Value = 0.1111111111111111
$ws.Range("H13").Value = 0.07142857142857142
$ws.Range("I13").Value = 0.07692307692307693
$ws.Range("J13").Value = 0.1111111111111111
$ws.Range("K13").Value = 0.07142857142857142
$ws.Range("L13").Value = 0.07692307692307693

# --- Recomputed metrics for row 14 ---
$ws.Range("G14").Value = 0.4
$ws.Range("H14").Value = 0.1
$ws.Range("I14").Value = 0.1176470588235294
$ws.Range("J14").Value = 0.4
$ws.Range("K14").Value = 0.1
$ws.Range("L14").Value = 0.1176470588235294
